$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.323.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.249.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.77"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.89%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.30"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.593.33"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.333.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.66"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.124.64"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0972"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.19"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.11"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "38.95"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +7.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.77%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.14"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.40"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0798"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.25"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.62"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.69"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.31%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0301"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.756.85"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.22"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.19"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.51%  "
